# Atualizacao de bases das ligas, do dia: 15-04-2024 as 22:35
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 248 (id 246) with corrected/expanded odds data ---
$ws.Range("B248").Value = 6775591
$ws.Range("E248").Value = 45395.625
$ws.Range("F248").Value = "Rakow Czestochowa"
$ws.Range("G248").Value = "Legia Warsaw"
$ws.Range("H248").Value = 1
$ws.Range("I248").Value = 1
$ws.Range("J248").Value = "D"
$ws.Range("K248").Value = 2.25
$ws.Range("L248").Value = 3.3
$ws.Range("M248").Value = 3.2
$ws.Range("N248").Value = 2.05
$ws.Range("O248").Value = 3.4
$ws.Range("P248").Value = 3.6
$ws.Range("Q248").Value = -0.5
$ws.Range("R248").Value = 2.025
$ws.Range("S248").Value = 1.825
$ws.Range("T248").Value = 2.5
$ws.Range("U248").Value = 2.025
$ws.Range("V248").Value = 1.825
$ws.Range("W248").Value = -1
$ws.Range("X248").Value = 2.4
$ws.Range("Y248").Value = -1
$ws.Range("Z248").Value = -1
$ws.Range("AA248").Value = 0.825
$ws.Range("AB248").Value = -1
$ws.Range("AC248").Value = 0.825

# --- Append new row 249 (id 247) ---
$ws.Range("A249").Value = 247
$ws.Range("B249").Value = 6775593
$ws.Range("C249").Value = "Poland Ekstraklasa"
$ws.Range("D249").Value = "Poland Ekstraklasa"
$ws.Range("E249").Value = 45396.3125
$ws.Range("F249").Value = "LKS Lodz"
$ws.Range("G249").Value = "Radomiak Radom"
$ws.Range("H249").Value = 3
$ws.Range("I249").Value = 2
$ws.Range("J249").Value = "H"
$ws.Range("K249").Value = 3.5
$ws.Range("L249").Value = 3.4
$ws.Range("M249").Value = 2.05
$ws.Range("N249").Value = 3.75
$ws.Range("O249").Value = 3.5
$ws.Range("P249").Value = 1.95
$ws.Range("Q249").Value = 0.5
$ws.Range("R249").Value = 1.825
$ws.Range("S249").Value = 2.025
$ws.Range("T249").Value = 2.5
$ws.Range("U249").Value = 1.95
$ws.Range("V249").Value = 1.9
$ws.Range("W249").Value = 2.75
$ws.Range("X249").Value = -1
$ws.Range("Y249").Value = -1
$ws.Range("Z249").Value = 0.825
$ws.Range("AA249").Value = -1
$ws.Range("AB249").Value = 0.95
$ws.Range("AC249").Value = -1

# --- Append new row 250 (id 248) ---
$ws.Range("A250").Value = 248
$ws.Range("B250").Value = 6775589
$ws.Range("C250").Value = "Poland Ekstraklasa"
$ws.Range("D250").Value = "Poland Ekstraklasa"
$ws.Range("E250").Value = 45396.41666666666
$ws.Range("F250").Value = "Jagiellonia Bialystok"
$ws.Range("G250").Value = "Cracovia Krakow"
$ws.Range("H250").Value = 1
$ws.Range("I250").Value = 3
$ws.Range("J250").Value = "A"
$ws.Range("K250").Value = 1.833
$ws.Range("L250").Value = 3.5
$ws.Range("M250").Value = 3.8
$ws.Range("N250").Value = 1.75
$ws.Range("O250").Value = 3.6
$ws.Range("P250").Value = 4
$ws.Range("Q250").Value = -0.5
$ws.Range("R250").Value = 1.8
$ws.Range("S250").Value = 2.05
$ws.Range("T250").Value = 2.75
$ws.Range("U250").Value = 1.975
$ws.Range("V250").Value = 1.875
$ws.Range("W250").Value = -1
$ws.Range("X250").Value = -1
$ws.Range("Y250").Value = 3
$ws.Range("Z250").Value = -1
$ws.Range("AA250").Value = 1.05
$ws.Range("AB250").Value = 0.9750000000000001
$ws.Range("AC250").Value = -1

# --- Append new row 251 (id 249) ---
$ws.Range("A251").Value = 249
$ws.Range("B251").Value = 6775588
$ws.Range("C251").Value = "Poland Ekstraklasa"
$ws.Range("D251").Value = "Poland Ekstraklasa"
$ws.Range("E251").Value = 45396.52083333334
$ws.Range("F251").Value = "Gornik Zabrze"
$ws.Range("G251").Value = "Slask Wroclaw"
$ws.Range("H251").Value = 2
$ws.Range("I251").Value = 0
$ws.Range("J251").Value = "H"
$ws.Range("K251").Value = 2.25
$ws.Range("L251").Value = 3.2
$ws.Range("M251").Value = 3
$ws.Range("N251").Value = 2.3
$ws.Range("O251").Value = 3.1
$ws.Range("P251").Value = 3
$ws.Range("Q251").Value = -0.25
$ws.Range("R251").Value = 2.05
$ws.Range("S251").Value = 1.8
$ws.Range("T251").Value = 2.25
$ws.Range("U251").Value = 2.05
$ws.Range("V251").Value = 1.8
$ws.Range("W251").Value = 1.3
$ws.Range("X251").Value = -1
$ws.Range("Y251").Value = -1
$ws.Range("Z251").Value = 1.05
$ws.Range("AA251").Value = -1
$ws.Range("AB251").Value = -0.5
$ws.Range("AC251").Value = 0.4

# Replicate row 248's cell formatting (bold/centered id column, date format column)
# onto the newly appended rows 249-251, matching how the source sheet is formatted.
$ws.Range("A248:AC248").Copy()
$ws.Range("A249:AC251").PasteSpecial(-4122)
$excel.CutCopyMode = 0
